# Applies the commit "end-to-end owernship, customer journey,enterprise products"
# to the active document.
#
# Four localized edits inside the experience-bullet table:
#   1. "...collaborated on API integration..." -> "...collaborated on APIs integration..."
#   2. "Customer Needs & Requirements Discovery" -> "Customer Journey & Requirements Discovery"
#   3. "Project Management & Agile Delivery: " -> "End-to-End Ownership: "
#   4. "...end-to-end product development..." -> "...end-to-end enterprise product development..."

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "API integration" -> "APIs integration"
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("collaborated on API", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Change 1: anchor text not found" }
$rng.Collapse(0)
$rng.InsertAfter("s")
$insertedS = $d.Range($rng.End - 1, $rng.End)
# force a run-split at the insertion boundary (toggle a format on/off so the
# new character doesn't silently merge back into its neighbours' run)
$insertedS.Bold = 1
$insertedS.Bold = 0

# ---------------------------------------------------------------------------
# Change 2: "Customer Needs & Requirements Discovery" -> "Customer Journey & Requirements Discovery"
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Needs ", $true, $false, $false, $false, $false, $true, 1, $false, "Journey ", 2)
if (-not $found) { throw "Change 2: anchor text not found" }
$rng2 = $d.Content
$found = $rng2.Find.Execute("Journey ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Change 2: replacement text not found" }
$trailingSpace = $d.Range($rng2.End - 1, $rng2.End)
$trailingSpace.Bold = 0
$trailingSpace.Bold = 1

# ---------------------------------------------------------------------------
# Change 3: "Project Management & Agile Delivery: " -> "End-to-End Ownership: "
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Project Management & Agile Delivery: ", $true, $false, $false, $false, $false, $true, 1, $false, "End-to-End Ownership: ", 2)
if (-not $found) { throw "Change 3: anchor text not found" }
$rng2 = $d.Content
$found = $rng2.Find.Execute("End-to-End Ownership: ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Change 3: replacement text not found" }
$trailingColon = $d.Range($rng2.End - 2, $rng2.End)
$trailingColon.Bold = 0
$trailingColon.Bold = 1

# ---------------------------------------------------------------------------
# Change 4: "end-to-end product development" -> "end-to-end enterprise product development"
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("end-to-end product development", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Change 4: anchor text not found" }
$prodStart = $rng.Start + 11   # length of "end-to-end "
$prodEnd = $prodStart + 7      # length of "product"
$prodRng = $d.Range($prodStart, $prodEnd)
$prodRng.Text = "enterprise product"

Write-Output "All four edits applied."
